# Auto-generated: apply scheduled market-data refresh to Sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 23362.076
$ws.Range("J51").Value = 9661.6
$ws.Range("L51").Value = 9661.6
$ws.Range("N51").Value = -10629.6
$ws.Range("H62").Value = 19233458
$ws.Range("I62").Value = 31252422
$ws.Range("K62").Value = 31252422
$ws.Range("M62").Value = -31251798
$ws.Range("H65").Value = 19233458
$ws.Range("I65").Value = 31252422
$ws.Range("K65").Value = 156262110
$ws.Range("M65").Value = -156258990
$ws.Range("H80").Value = 1140.75
$ws.Range("I80").Value = 341.16666
$ws.Range("K80").Value = 1023.49998
$ws.Range("M80").Value = -25.49997999999994
$ws.Range("H83").Value = 1140.75
$ws.Range("I83").Value = 341.16666
$ws.Range("K83").Value = 3070.49994
$ws.Range("M83").Value = 1921.50006
$ws.Range("H100").Value = 5846.643
$ws.Range("I100").Value = 1498.5
$ws.Range("J100").Value = 7585.9
$ws.Range("K100").Value = 1498.5
$ws.Range("L100").Value = 7585.9
$ws.Range("M100").Value = -957.5
$ws.Range("N100").Value = -8667.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3954.85
$ws.Range("I32").Value = 3651.8438
$ws.Range("K32").Value = 3651.8438
$ws.Range("M32").Value = -3364.8438
$ws.Range("H45").Value = 2767
$ws.Range("I45").Value = 1159.625
$ws.Range("J45").Value = 5981.75
$ws.Range("K45").Value = 1159.625
$ws.Range("L45").Value = 5981.75
$ws.Range("M45").Value = -782.625
$ws.Range("N45").Value = -6735.75
$ws.Range("H61").Value = 10640.333
$ws.Range("I61").Value = 6010.5557
$ws.Range("K61").Value = 6010.5557
$ws.Range("M61").Value = -5798.5557
$ws.Range("H63").Value = 2123
$ws.Range("I63").Value = 2123
$ws.Range("K63").Value = 2123
$ws.Range("M63").Value = -1437
$ws.Range("H66").Value = 2123
$ws.Range("I66").Value = 2123
$ws.Range("K66").Value = 10615
$ws.Range("M66").Value = -7183
$ws.Range("H74").Value = 2294
$ws.Range("I74").Value = 1356.7
$ws.Range("J74").Value = 3335.4443
$ws.Range("K74").Value = 1356.7
$ws.Range("L74").Value = 3335.4443
$ws.Range("M74").Value = -482.7
$ws.Range("N74").Value = -5083.4443
$ws.Range("H77").Value = 2294
$ws.Range("I77").Value = 1356.7
$ws.Range("J77").Value = 3335.4443
$ws.Range("K77").Value = 6783.5
$ws.Range("L77").Value = 16677.2215
$ws.Range("M77").Value = -2415.5
$ws.Range("N77").Value = -25413.2215
$ws.Range("H132").Value = 4081.6667
$ws.Range("I132").Value = 3478.1
$ws.Range("J132").Value = 7099.5
$ws.Range("K132").Value = 10434.3
$ws.Range("L132").Value = 21298.5
$ws.Range("M132").Value = -7904.299999999999
$ws.Range("N132").Value = -26358.5
$ws.Range("H136").Value = 10640.333
$ws.Range("I136").Value = 6010.5557
$ws.Range("K136").Value = 18031.6671
$ws.Range("M136").Value = -15481.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1910.7778
$ws.Range("I105").Value = 1937.9412
$ws.Range("J105").Value = 1449
$ws.Range("K105").Value = 1937.9412
$ws.Range("L105").Value = 1449
$ws.Range("M105").Value = -190.9412
$ws.Range("N105").Value = -4943
$ws.Range("H107").Value = 3287.6
$ws.Range("J107").Value = 2776.7144
$ws.Range("L107").Value = 2776.7144
$ws.Range("N107").Value = -6616.7144
$ws.Range("H108").Value = 212990
$ws.Range("J108").Value = 212990
$ws.Range("L108").Value = 212990
$ws.Range("N108").Value = -220670
$ws.Range("H134").Value = 5834.263
$ws.Range("I134").Value = 5697.3125
$ws.Range("K134").Value = 17091.9375
$ws.Range("M134").Value = -14556.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 114.875
$ws.Range("J7").Value = 222
$ws.Range("L7").Value = 222
$ws.Range("N7").Value = -448
$ws.Range("H16").Value = 4261.5
$ws.Range("I16").Value = 2894.75
$ws.Range("J16").Value = 6995
$ws.Range("K16").Value = 2894.75
$ws.Range("L16").Value = 6995
$ws.Range("M16").Value = -2607.75
$ws.Range("N16").Value = -7569
$ws.Range("H58").Value = 2957.5
$ws.Range("I58").Value = 1653.909
$ws.Range("K58").Value = 1653.909
$ws.Range("M58").Value = -1450.909
$ws.Range("H103").Value = 24438.8
$ws.Range("I103").Value = 12100
$ws.Range("K103").Value = 12100
$ws.Range("M103").Value = -10928
$ws.Range("H113").Value = 4261.5
$ws.Range("I113").Value = 2894.75
$ws.Range("J113").Value = 6995
$ws.Range("K113").Value = 2894.75
$ws.Range("L113").Value = 6995
$ws.Range("M113").Value = -724.75
$ws.Range("N113").Value = -11335
$ws.Range("H132").Value = 3246.3157
$ws.Range("I132").Value = 3088.3333
$ws.Range("K132").Value = 9264.999899999999
$ws.Range("M132").Value = -6734.999899999999
$ws.Range("H134").Value = 4973.484
$ws.Range("I134").Value = 3970.3215
$ws.Range("K134").Value = 11910.9645
$ws.Range("M134").Value = -9375.9645
$ws.Range("H136").Value = 2957.5
$ws.Range("I136").Value = 1653.909
$ws.Range("K136").Value = 4961.727000000001
$ws.Range("M136").Value = -2411.727000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2317.8
$ws.Range("I7").Value = 120
$ws.Range("K7").Value = 360
$ws.Range("M7").Value = -248
$ws.Range("H33").Value = 1422.625
$ws.Range("J33").Value = 199
$ws.Range("L33").Value = 1194
$ws.Range("N33").Value = -1760
$ws.Range("H40").Value = 2447.5625
$ws.Range("J40").Value = 4560.8
$ws.Range("L40").Value = 18243.2
$ws.Range("N40").Value = -18381.2
$ws.Range("H51").Value = 1466.6666
$ws.Range("I51").Value = 1173.125
$ws.Range("J51").Value = 1802.1428
$ws.Range("K51").Value = 3519.375
$ws.Range("L51").Value = 5406.428400000001
$ws.Range("M51").Value = -3059.375
$ws.Range("N51").Value = -6326.428400000001
$ws.Range("H122").Value = 3143.1428
$ws.Range("I122").Value = 398.6
$ws.Range("K122").Value = 3587.4
$ws.Range("M122").Value = -1137.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 6999
$ws.Range("I5").Value = 4999
$ws.Range("J5").Value = 9999
$ws.Range("K5").Value = 4999
$ws.Range("L5").Value = 9999
$ws.Range("M5").Value = -4887
$ws.Range("N5").Value = -10223
$ws.Range("H10").Value = 81285.28999999999
$ws.Range("I10").Value = 168333.33
$ws.Range("J10").Value = 15999.25
$ws.Range("K10").Value = 168333.33
$ws.Range("L10").Value = 15999.25
$ws.Range("M10").Value = -168164.33
$ws.Range("N10").Value = -16337.25
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("M12").Value = 40
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H55").Value = 12999.667
$ws.Range("I55").Value = 9000
$ws.Range("J55").Value = 14999.5
$ws.Range("K55").Value = 9000
$ws.Range("L55").Value = 14999.5
$ws.Range("M55").Value = -8673
$ws.Range("N55").Value = -15653.5
$ws.Range("H113").Value = 3935
$ws.Range("I113").Value = 3722
$ws.Range("K113").Value = 3722
$ws.Range("M113").Value = -1552
$ws.Range("H132").Value = 2940.4546
$ws.Range("I132").Value = 2934.5
$ws.Range("K132").Value = 8803.5
$ws.Range("M132").Value = -6273.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 49999
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H46").Value = 1923.0625
$ws.Range("I46").Value = 1052.7778
$ws.Range("K46").Value = 1052.7778
$ws.Range("M46").Value = -864.7778000000001
$ws.Range("H48").Value = 49999
$ws.Range("J48").Value = 49999
$ws.Range("L48").Value = 49999
$ws.Range("N48").Value = -51321
$ws.Range("H132").Value = 4457.1665
$ws.Range("I132").Value = 3438
$ws.Range("J132").Value = 5185.143
$ws.Range("K132").Value = 10314
$ws.Range("L132").Value = 15555.429
$ws.Range("M132").Value = -7784
$ws.Range("N132").Value = -20615.429
$ws.Range("H136").Value = 2278.3865
$ws.Range("I136").Value = 1157.5834
$ws.Range("K136").Value = 3472.7502
$ws.Range("M136").Value = -922.7501999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3055.125
$ws.Range("I96").Value = 1999.8
$ws.Range("J96").Value = 4814
$ws.Range("K96").Value = 1999.8
$ws.Range("L96").Value = 4814
$ws.Range("M96").Value = -626.8
$ws.Range("N96").Value = -7560
$ws.Range("H100").Value = 1110.8572
$ws.Range("I100").Value = 461.5
$ws.Range("K100").Value = 923
$ws.Range("M100").Value = -382
$ws.Range("H122").Value = 2311.575
$ws.Range("I122").Value = 1456.7812
$ws.Range("J122").Value = 5730.75
$ws.Range("K122").Value = 4370.3436
$ws.Range("L122").Value = 17192.25
$ws.Range("M122").Value = -1920.3436
$ws.Range("N122").Value = -22092.25
$ws.Range("H126").Value = 2527.5
$ws.Range("I126").Value = 2291.25
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6873.75
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -4403.75
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 1362.75
$ws.Range("I132").Value = 1308.2307
$ws.Range("K132").Value = 3924.6921
$ws.Range("M132").Value = -1394.6921
$ws.Range("H136").Value = 10180.667
$ws.Range("I136").Value = 11766.85
$ws.Range("K136").Value = 35300.55
$ws.Range("M136").Value = -32750.55

